$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.448.70'
$ws.Range('E2').Value = '  -0.02%  '

$ws.Range('D3').Value = '1.810.19'
$ws.Range('E3').Value = '  +0.39%  '

$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.70%  '

$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  -0.65%  '

$ws.Range('D6').Value = '305.23'
$ws.Range('E6').Value = '  -1.01%  '

$ws.Range('D7').Value = '0.4515'
$ws.Range('E7').Value = '  -0.77%  '

$ws.Range('D8').Value = '0.3596'
$ws.Range('E8').Value = '  -1.58%  '

$ws.Range('D9').Value = '46.22'
$ws.Range('E9').Value = '  +2.64%  '

$ws.Range('D10').Value = '0.07058'
$ws.Range('E10').Value = '  -0.85%  '

$ws.Range('D11').Value = '0.8905'
$ws.Range('E11').Value = '  +1.51%  '

$ws.Range('D12').Value = '0.07782'
$ws.Range('E12').Value = '  +0.61%  '

$ws.Range('D13').Value = '19.32'
$ws.Range('E13').Value = '  +0.04%  '

$ws.Range('D14').Value = '1.777.14'
$ws.Range('E14').Value = '  -1.79%  '

$ws.Range('D15').Value = '5.270'
$ws.Range('E15').Value = '  +0.06%  '

$ws.Range('D16').Value = '6.301'
$ws.Range('E16').Value = '  -0.75%  '

$ws.Range('D17').Value = '84.86'
$ws.Range('E17').Value = '  -1.25%  '

$ws.Range('D18').Value = '1.004'
$ws.Range('E18').Value = '  -0.65%  '

$ws.Range('D19').Value = '0.000008518'
$ws.Range('E19').Value = '  -0.60%  '

$ws.Range('E20').Value = '  -0.58%  '

$ws.Range('D21').Value = '26.482.79'

$ws.Range('D22').Value = '14.19'
$ws.Range('E22').Value = '  -0.30%  '

$ws.Range('D23').Value = '4.953'
$ws.Range('E23').Value = '  -0.35%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '10.54'
$ws.Range('E24').Value = '  +1.22%  '

$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.015.59'
$ws.Range('E25').Value = '  -1.38%  '

$ws.Range('D26').Value = '1.961'
$ws.Range('E26').Value = '  -1.35%  '

$ws.Range('D27').Value = '150.80'
$ws.Range('E27').Value = '  +0.04%  '

$ws.Range('D28').Value = '17.77'
$ws.Range('E28').Value = '  -0.82%  '

$ws.Range('D29').Value = '2.045'
$ws.Range('E29').Value = '  +2.02%  '

$ws.Range('D30').Value = '112.23'
$ws.Range('E30').Value = '  -0.18%  '

$ws.Range('D31').Value = '4.825'
$ws.Range('E31').Value = '  +0.25%  '

$ws.Range('D32').Value = '0.08686'

$ws.Range('D33').Value = '3.132'
$ws.Range('E33').Value = '  +2.69%  '

$ws.Range('D34').Value = '0.7454'
$ws.Range('E34').Value = '  +2.38%  '

$ws.Range('D35').Value = '2.740'
$ws.Range('E35').Value = '  +7.34%  '

$ws.Range('D36').Value = '4.420'
$ws.Range('E36').Value = '  -0.21%  '

$ws.Range('D37').Value = '1.109'
$ws.Range('E37').Value = '  -0.31%  '

$ws.Range('D38').Value = '1.065'
$ws.Range('E38').Value = '  -1.23%  '

$ws.Range('D39').Value = '0.01928'
$ws.Range('E39').Value = '  -0.06%  '

$ws.Range('D40').Value = '2.900'
$ws.Range('E40').Value = '  +0.74%  '

$ws.Range('D41').Value = '0.05090'
$ws.Range('E41').Value = '  -0.08%  '

$ws.Range('D42').Value = '0.5079'
$ws.Range('E42').Value = '  +1.61%  '

$ws.Range('D43').Value = '6.730'
$ws.Range('E43').Value = '  -2.94%  '

$ws.Range('D44').Value = '0.1506'
$ws.Range('E44').Value = '  -3.60%  '

$ws.Range('D45').Value = '8.036'
$ws.Range('E45').Value = '  -0.81%  '

$ws.Range('D46').Value = '0.4719'

$ws.Range('D47').Value = '1.002'
$ws.Range('E47').Value = '  -0.76%  '

$ws.Range('D48').Value = '9.989'
$ws.Range('E48').Value = '  +1.00%  '

$ws.Range('D49').Value = '100.20'
$ws.Range('E49').Value = '  -1.28%  '

$ws.Range('D50').Value = '1.575'
$ws.Range('E50').Value = '  -0.84%  '

$ws.Range('D51').Value = '0.05980'
$ws.Range('E51').Value = '  -0.09%  '
